$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure D and E columns are formatted as text so numeric-looking strings
# (e.g. "553.73", "64.753.28") are preserved exactly as strings, matching
# the original inlineStr cell type.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "64.753.28"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").Value = "3.349.89"
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "553.73"
$ws.Range("E5").Value = "  -0.13%  "

$ws.Range("D6").Value = "173.60"
$ws.Range("E6").Value = "  -0.81%  "

$ws.Range("E7").Value = "  +2.07%  "

$ws.Range("D8").Value = "3.335.96"
$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  -0.07%  "

$ws.Range("E10").Value = "  +6.88%  "

$ws.Range("D11").Value = "0.636"
$ws.Range("E11").Value = "  +1.61%  "

$ws.Range("D12").Value = "53.47"
$ws.Range("E12").Value = "  -2.02%  "

$ws.Range("D13").Value = "0.0000280"
$ws.Range("E13").Value = "  +3.43%  "

$ws.Range("E14").Value = "  +0.78%  "

$ws.Range("D15").Value = "3.883.46"
$ws.Range("E15").Value = "  -0.20%  "

$ws.Range("E16").Value = "  +2.53%  "

$ws.Range("E17").Value = "  -0.40%  "

$ws.Range("D18").Value = "3.348.49"
$ws.Range("E18").Value = "  -0.13%  "

$ws.Range("D19").Value = "64.725.15"
$ws.Range("E19").Value = "  +0.57%  "

$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("E21").Value = "  +1.09%  "

$ws.Range("D22").Value = "446.99"
$ws.Range("E22").Value = "  +2.67%  "

$ws.Range("D23").Value = "4.92"
$ws.Range("E23").Value = "  -0.84%  "

$ws.Range("D24").Value = "4.06"
$ws.Range("E24").Value = "  -0.34%  "

$ws.Range("D25").Value = "86.79"
$ws.Range("E25").Value = "  +2.96%  "

$ws.Range("D26").Value = "13.64"
$ws.Range("E26").Value = "  +1.68%  "

$ws.Range("E27").Value = "  +1.67%  "

$ws.Range("D28").Value = "10.69"
$ws.Range("E28").Value = "  -0.34%  "

$ws.Range("D29").Value = "8.62"
$ws.Range("E29").Value = "  -0.90%  "

$ws.Range("E30").Value = "  +3.97%  "

$ws.Range("E31").Value = "  -1.53%  "

$ws.Range("D32").Value = "62.92"
$ws.Range("E32").Value = "  +7.82%  "

$ws.Range("D33").Value = "11.42"
$ws.Range("E33").Value = "  -0.20%  "

$ws.Range("D34").Value = "574.66"
$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("E35").Value = "  -0.26%  "

$ws.Range("E36").Value = "  +0.02%  "

$ws.Range("D37").Value = "3.62"
$ws.Range("E37").Value = "  +2.92%  "

$ws.Range("E38").Value = "  -0.15%  "

$ws.Range("D39").Value = "35.53"
$ws.Range("E39").Value = "  -0.11%  "

$ws.Range("D40").Value = "0.0₃0740"
$ws.Range("E40").Value = "  -1.44%  "

$ws.Range("D41").Value = "0.369"
$ws.Range("E41").Value = "  +0.90%  "

$ws.Range("D42").Value = "3.077.42"
$ws.Range("E42").Value = "  -0.64%  "

$ws.Range("E43").Value = "  +1.93%  "

$ws.Range("E44").Value = "  -2.10%  "

$ws.Range("E45").Value = "  +3.64%  "

$ws.Range("E48").Value = "  -0.03%  "

$ws.Range("D49").Value = "140.69"
$ws.Range("E49").Value = "  +4.11%  "

$ws.Range("D50").Value = "2.52"
$ws.Range("E50").Value = "  -2.60%  "

$ws.Range("D51").Value = "8.25"
$ws.Range("E51").Value = "  -0.08%  "

# Rows 46 and 47: ApeXProtocol and Fetch.AI swap places
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "2.45"
$ws.Range("E46").Value = "  -0.27%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "3.17"
$ws.Range("E47").Value = "  -1.24%  "
